$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1 (Title 1): "A Table, with a caption"
# Consolidate run pairs ("word" + " ") into a single run with a trailing space.
$tr1 = $s.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 2).Text = "A "
$tr1.Characters(3, 7).Text = "Table, "
$tr1.Characters(10, 5).Text = "with "
$tr1.Characters(15, 2).Text = "a "

# Shape 3 (TextBox 3): "Demonstration of simple table syntax, with alignment"
$tr3 = $s.Shapes.Item(3).TextFrame.TextRange
$tr3.Characters(1, 14).Text = "Demonstration "
$tr3.Characters(15, 3).Text = "of "
$tr3.Characters(18, 7).Text = "simple "
$tr3.Characters(25, 6).Text = "table "
$tr3.Characters(31, 8).Text = "syntax, "
$tr3.Characters(39, 5).Text = "with "
